# Lactose operon model - rename parameters to the new short naming scheme
# (alpha_X -> a_N, beta_X -> b_N, gamma_X -> g_N, Gamma_0 -> G_0) and fix a
# couple of subscript typos in the reaction-term formulas, then add the new
# K_6 parameter (and carry K_5 down into its own row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Parameter names + reaction-term formulas (columns D/E, rows 2-6) ---
$ws.Range("D2").Value = "a_1"
$ws.Range("E2").Value = "a_1*((1+K_1*A)/(K_0+K_1*A))+G_0-g_1*M"

$ws.Range("D3").Value = "a_2"
$ws.Range("E3").Value = "a_2*M-g_2*B"

$ws.Range("D4").Value = "a_3"
$ws.Range("E4").Value = "a_3*((P*L)/(K_2+L))-b_1*((P*L)/(K_3+L))-b_2*((B*L)/(K_4+L))-g_3*L"

$ws.Range("D5").Value = "a_4"
$ws.Range("E5").Value = "a_4*((B*L)/(K_5+L))-b_3*((B*A)/(K_6+A))-g_4*A"

$ws.Range("D6").Value = "a_5"
$ws.Range("E6").Value = "a_5*M-g_5*P"

# --- Parameter list (column D, rows 7-20) ---
$ws.Range("D7").Value  = "b_1"
$ws.Range("D8").Value  = "b_2"
$ws.Range("D9").Value  = "b_3"
$ws.Range("D10").Value = "g_1"
$ws.Range("D11").Value = "g_2"
$ws.Range("D12").Value = "g_3"
$ws.Range("D13").Value = "g_4"
$ws.Range("D14").Value = "g_5"
$ws.Range("D15").Value = "G_0"
$ws.Range("D16").Value = "K_0"
$ws.Range("D17").Value = "K_1"
$ws.Range("D18").Value = "K_2"
$ws.Range("D19").Value = "K_3"
$ws.Range("D20").Value = "K_4"

# --- New rows for K_5 and the newly introduced K_6 ---
$ws.Range("D21").Value = "K_5"
$ws.Range("D22").Value = "K_6"

# Move the active selection (the author was last looking at the a_4 formula)
$ws.Range("E5").Select() | Out-Null
